# Apply the weather-data refresh update to cfb_weather.xlsx
$wb = $excel.ActiveWorkbook

$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# --- Update the run Timestamp column (shared across every row on FBS) ---
$wsFBS.Range("AK2:AK59").Value = "2024-11-21T10:01:54.388753"

# --- FBS sheet numeric / spread-move updates ---
# Row 18: Iowa @ Maryland
$wsFBS.Range("AB18").Value = 6.5
$wsFBS.Range("AF18").Value = -0.5

# Row 34: Penn State @ Minnesota
$wsFBS.Range("Q34").Value = "SSE"

# Row 35: San Diego State @ Utah State
$wsFBS.Range("Q35").Value = "E"

# Row 37: Stanford @ California
$wsFBS.Range("AB37").Value = -14
$wsFBS.Range("AF37").Value = 0

# Row 39: Pittsburgh @ Louisville
$wsFBS.Range("AB39").Value = -8.5
$wsFBS.Range("AF39").Value = 0

# Row 46: Marshall @ Old Dominion
$wsFBS.Range("AB46").Value = -1.5
$wsFBS.Range("AF46").Value = -1

# Row 47: Iowa State @ Utah
$wsFBS.Range("Q47").Value = "SSE"
$wsFBS.Range("Y47").Value = 41.5
$wsFBS.Range("AE47").Value = 0

# Row 48: Texas A&M @ Auburn
$wsFBS.Range("Q48").Value = "SSE"

# Row 53: USC @ UCLA
$wsFBS.Range("Q53").Value = "NNW"

# Row 56: Western Michigan @ Central Michigan
$wsFBS.Range("Q56").Value = "NE"

# --- Other sheet wind_dir_fg updates ---
# Row 30: Butler vs Presbyterian
$wsOther.Range("S30").Value = "ENE"

# Row 42: South Dakota State vs Missouri State
$wsOther.Range("S42").Value = "N"

# Row 43: Cal Poly vs Weber State
$wsOther.Range("S43").Value = "NE"
